# Kostya - update Browser users
# Adds 10 new "AUtestuser" rows (with their Email/Role/Password columns and a
# new "H" (answer) column) into the Users sheet, right after the existing
# AU annotation-user block (old row 17) and before the blank separator row
# that used to be row 18 (now row 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$xlPasteFormats = -4122
$xlNone = -4142
$xlEdgeTop = 8
$xlEdgeBottom = 9

# ---------------------------------------------------------------------------
# 1. Insert 10 new blank rows at 18..27 (pushes the old rows 18-24 down to
#    28-34, matching the target layout).
# ---------------------------------------------------------------------------
$ws.Range("A18:A27").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Copy formatting for the new rows from the existing "plain" data rows so
#    the new cells pick up the same borders/fonts as the rest of the table.
#    Row 9 (A:s6, B:s6, C:s6, D:s6, E:s6, F:s6, G:s7-hyperlink) is a good
#    template for A:G. Column A in the new rows additionally needs the
#    "s8" look (plain font but explicitly applied) like rows 7/8 use.
# ---------------------------------------------------------------------------
for ($r = 18; $r -le 27; $r++) {
    $ws.Range("A9:G9").Copy()
    $ws.Range("A$r" + ":G$r").PasteSpecial($xlPasteFormats)
    $ws.Range("A7").Copy()
    $ws.Range("A$r").PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the values for the 10 new AU test users.
# ---------------------------------------------------------------------------
$names  = @("AUtestuser1","AUtestuser2","AUtestuser3","AUtestuser4","AUtestuser5","AUtestuser6","AUtestuser7","AUtestuser8","AUtestuser9","AUtestuser10")
$roles  = @("Admin","Admin","","","","","","","","")
$emails = @("AUtestuser1@mailinator.com","AUtestuser2@mailinator.com","AUtestuser3@mailinator.com","AUtestuser4@mailinator.com","AUtestuser5@mailinator.com","AUtestuser6@mailinator.com","AUtestuser7@mailinator.com","AUtestuser8@mailinator.com","AUtestuser9@mailinator.com","AUtestuser10@mailinator.com")

for ($i = 0; $i -lt 10; $i++) {
    $r = 18 + $i
    $ws.Range("A$r").Value = $names[$i]
    $ws.Range("B$r").Value = "Password1"
    if ($roles[$i] -ne "") {
        $ws.Range("D$r").Value = $roles[$i]
    }
    $ws.Range("G$r").Value = $emails[$i]
    $ws.Range("H$r").Value = "thomsonreuters"
}

# ---------------------------------------------------------------------------
# 4. New header cell H1 ("answer") - same bold/fill header look as A1:G1 but
#    with a border that only has left/right (no top/bottom).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial($xlPasteFormats)
$ws.Range("H1").Borders.Item($xlEdgeTop).LineStyle = $xlNone
$ws.Range("H1").Borders.Item($xlEdgeBottom).LineStyle = $xlNone
$ws.Range("H1").Value = "answer"
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. G7's hyperlink-styled cell loses its explicit "applyFont" (cosmetic
#    restyle in the source workbook) - re-apply the same visual hyperlink
#    format that G9 already carries (font/border without forcing applyFont).
# ---------------------------------------------------------------------------
$ws.Range("G9").Copy()
$ws.Range("G7").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Rebuild every hyperlink on the sheet (row-insert does not shift the
#    existing hyperlink anchors in this engine, so clear + re-add them all
#    at their correct, final addresses).
# ---------------------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:anzuser1@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:anzuser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:anzuser3@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:anzuser4@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G7"), "mailto:anztestuser1@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G8"), "mailto:anztestuser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G10"), "mailto:anztestuser4@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G11"), "mailto:anzemployment@mailinator.com")

$ws.Hyperlinks.Add($ws.Range("G18"), "mailto:AUtestuser1@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G19"), "mailto:AUtestuser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G20"), "mailto:AUtestuser3@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G21"), "mailto:AUtestuser4@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G22"), "mailto:AUtestuser5@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G23"), "mailto:AUtestuser6@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G24"), "mailto:AUtestuser7@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G25"), "mailto:AUtestuser8@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G26"), "mailto:AUtestuser9@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G27"), "mailto:AUtestuser10@mailinator.com")

$ws.Hyperlinks.Add($ws.Range("G32"), "mailto:shareannotationuser1@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G33"), "mailto:myShareAnnotationUser@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G30"), "mailto:auannotationuser2@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("G31"), "mailto:auannotationuser3@mailinator.com")

# ---------------------------------------------------------------------------
# 7. Selection cursor moved to K24 in the source edit.
# ---------------------------------------------------------------------------
$ws.Range("K24").Select()
